$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new table (mirrors A1:D1)
$ws.Range("F1").Value = "CUDA"
$ws.Range("G1").Value = 512
$ws.Range("H1").Value = 1024
$ws.Range("I1").Value = 2048

# Row 2: "Global" speed-up ratios vs. CPU (row 3)
$ws.Range("F2").Value = "Global"
$ws.Range("G2").Formula = "=B2/B3"
$ws.Range("H2").Formula = "=C2/C3"
$ws.Range("I2").Formula = "=D2/D3"

# Row 3: "Shared" speed-up ratios vs. CPU (row 4)
$ws.Range("F3").Value = "Shared"
$ws.Range("G3").Formula = "=B2/B4"
$ws.Range("H3").Formula = "=C2/C4"
$ws.Range("I3").Formula = "=D2/D4"

$ws.Range("F4").Select()
